# Outstandings.xlsx - "Add files via upload" edit
#
# Summary of the change (reconstructed from the OOXML diff):
#  - The worksheet "Purchase 22-23" contained two copies of the same
#    purchase ledger: rows 1-19 (table #1) and rows 22-37 (table #2, an
#    exact duplicate of rows 1-16 plus the two extra rows 18-19).
#  - Table #2 (rows 22-37, including its own header row 22) is deleted.
#  - A new ledger entry is inserted as row 17 in table #1 (date
#    11-Jul-2023, invoice "30/23-23", client "Nerul", amount 31907,
#    running balance formula F16+E17), pushing the former rows 18-19
#    down to rows 19-20.
#  - The sheet view/selection & dimension shrink to match the now much
#    smaller used range.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Purchase 22-23")

# 1) Remove the duplicated second table (old header row 22 + data rows
#    23-37) entirely.
$ws.Rows("22:37").Delete() | Out-Null

# 2) Row 16's running-balance cell (F16) picks up the plain "s=9" style
#    (previously "s=25" - the same visual style, just a different
#    cellXf entry) instead of the accent style used for blank summary
#    cells. Re-use the format that's already sitting in A16 (style 9).
$ws.Range("A16").Copy() | Out-Null
$ws.Range("F16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# 3) Insert a brand-new blank row above the old row 18 - this shifts the
#    old rows 18 & 19 down to 19 & 20 (their values/relative formulas
#    auto-adjust, e.g. "=E18+E19" becomes "=E19+E20").
$ws.Rows("17:17").Insert() | Out-Null

# The freshly inserted row 17 has no formatting yet; row 19 (the old
# row 18) already carries the exact formats we want for every column
# of row 17 (s=9, s=24, s=23, s=23, s=23, s=25), so copy them down.
$ws.Range("A19:F19").Copy() | Out-Null
$ws.Range("A17:F17").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Rows("17:17").RowHeight = $ws.Rows("19:19").RowHeight

# 4) Populate the new ledger row. Column D ("Vendor") re-uses the same
#    vendor name already used by every other row in this table.
$ws.Range("B17").Value = 45118
$ws.Range("C17").Value = "30/23-23"
$ws.Range("D17").Value = $ws.Range("D16").Value2
$ws.Range("E17").Value = 31907
$ws.Range("F17").Formula = "=F16+E17"

# 5) Reflect the new (much smaller) used range & cursor position.
$ws.Range("A1:F20").Select() | Out-Null
$ws.Range("G23").Select() | Out-Null

$excel.CutCopyMode = $false
